# study1_result2.xlsx - "Add files via upload" edit
#
# The old layout had a throwaway numeric index column in column A
# (0..6) and a row of plain sequential numbers (0..8) in row 1 across
# columns B:J. The new layout drops that index column entirely (so
# everything shifts one column to the left) and replaces the numeric
# row-1 header with real column names in A1:I1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old index column (A). This shifts B:J left to A:I,
# carrying cell contents AND formatting (so the header style that was
# on B1:J1 ends up on A1:I1, and the un-styled data in B2:B8 ends up
# un-styled in A2:A8, matching the target).
$ws.Columns.Item(1).Delete()

# Overwrite the (now former-B) header row with the real column names.
$headers = @("subsets", "method", "F_val", "dfb", "dfw", "p_val", "partial_eta2", "cohens_f", "post-hoc mean chain")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$wb.Save()
